$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Fix the header labels: hyphens -> underscores
$ws.Range("B1").Value = "aTPJ_R_P_F_C"
$ws.Range("C1").Value = "aTPJ_R_P_F_S"
$ws.Range("D1").Value = "pTPJ_R_P_F_C"
$ws.Range("E1").Value = "pTPJ_R_P_F_S"

# 2. Insert a missing subject row (1303) right after row 25 (subject 1302),
#    shifting the existing rows 26-44 down to 27-45.
$insertRow = $ws.Rows.Item(26)
$insertRow.Insert(-4121) # xlShiftDown

$ws.Range("A26").Value = 1303
$ws.Range("B26").Value = 0.0829808524
$ws.Range("C26").Value = -0.02110848941999999
$ws.Range("D26").Value = -0.00618960030999996
$ws.Range("E26").Value = -0.16397736683

# 3. Append the missing subject row (3220) at the new end of the data (row 46).
$ws.Range("A46").Value = 3220
$ws.Range("B46").Value = 0.07862091610000001
$ws.Range("C46").Value = -0.08206761560000002
$ws.Range("D46").Value = 0.3183748924300001
$ws.Range("E46").Value = -0.04529498433000001
